$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row (row 16) of abbreviated header names.
# Cells are set in an order that matches the shared-string table
# growth order recorded in the target workbook (new unique strings
# are appended to xl/sharedStrings.xml in first-use order).

# K16 / L16 reuse existing shared strings ("Brand" / "Model").
$ws.Range("K16").Value = "Brand"
$ws.Range("L16").Value = "Model"

# New shared strings, in the exact order they first appear.
$ws.Range("A16").Value = "sno"
$ws.Range("B16").Value = "SeatNo"
$ws.Range("C16").Value = "EmpID"
$ws.Range("D16").Value = "EmpNm"
$ws.Range("F16").Value = "fb"
$ws.Range("G16").Value = "status"
$ws.Range("H16").Value = "atn"
$ws.Range("J16").Value = "Asstyp"
$ws.Range("I16").Value = "HostNm"
$ws.Range("E16").Value = "Loc"
$ws.Range("M16").Value = "SrlNo"
$ws.Range("N16").Value = "AssDev"
$ws.Range("O16").Value = "ADSNo"
$ws.Range("P16").Value = "HDD"
$ws.Range("Q16").Value = "Mem"
$ws.Range("R16").Value = "Proc"
$ws.Range("S16").Value = "OS"
$ws.Range("T16").Value = "PurchOn"
$ws.Range("U16").Value = "Inv"
$ws.Range("V16").Value = "Vend"
$ws.Range("W16").Value = "Wrty"
$ws.Range("X16").Value = "DoI"
$ws.Range("Y16").Value = "Cno"
$ws.Range("Z16").Value = "Email"
$ws.Range("AA16").Value = "Rmks"

# Update the active selection to A16, matching the saved view state.
$ws.Range("A16").Select()
